# Apply the "cryptos list" refresh: updated prices / 1h volume deltas,
# plus two pairs of rows (26/27 and 49/50) whose ranking order swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '58.628.37'
$ws.Range("E2").Value = '  -1.05%  '
# Row 3
$ws.Range("D3").Value = '2.627.90'
$ws.Range("E3").Value = '  -0.43%  '
# Row 4
$ws.Range("E4").Value = '  -0.02%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '520.28'
$ws.Range("E5").Value = '  +0.89%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.58'
$ws.Range("E6").Value = '  -3.26%  '
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.39%  '
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.575'
$ws.Range("E8").Value = '  -0.52%  '
# Row 9
$ws.Range("D9").Value = '2.638.04'
$ws.Range("E9").Value = '  -1.14%  '
# Row 10
$ws.Range("E10").Value = '  -4.28%  '
# Row 11
$ws.Range("E11").Value = '  -2.97%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.333'
$ws.Range("E12").Value = '  -2.36%  '
# Row 13
$ws.Range("E13").Value = '  -0.73%  '
# Row 14
$ws.Range("D14").Value = '3.089.33'
$ws.Range("E14").Value = '  -0.49%  '
# Row 15
$ws.Range("D15").Value = '58.653.37'
$ws.Range("E15").Value = '  -0.73%  '
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.73'
# Row 17
$ws.Range("E17").Value = '  -2.91%  '
# Row 18
$ws.Range("D18").Value = '2.633.93'
$ws.Range("E18").Value = '  -0.87%  '
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '346.31'
$ws.Range("E19").Value = '  -0.46%  '
# Row 20
$ws.Range("E20").Value = '  -4.20%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.19'
$ws.Range("E21").Value = '  -3.92%  '
# Row 22
$ws.Range("E22").Value = '  -2.60%  '
# Row 23
$ws.Range("E23").Value = '  +0.03%  '
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '61.69'
$ws.Range("E24").Value = '  +0.77%  '
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.413'
$ws.Range("E25").Value = '  -3.19%  '
# Row 26
$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.74%  '
# Row 27
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.163'
$ws.Range("E27").Value = '  +0.48%  '
# Row 28
$ws.Range("D28").Value = '0.0₃0799'
$ws.Range("E28").Value = '  -4.63%  '
# Row 29
$ws.Range("E29").Value = '  -1.96%  '
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.21'
$ws.Range("E31").Value = '  -3.24%  '
# Row 32
$ws.Range("E32").Value = '  -0.12%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '149.19'
$ws.Range("E34").Value = '  +0.02%  '
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.978'
$ws.Range("E35").Value = '  -6.59%  '
# Row 36
$ws.Range("E36").Value = '  -3.51%  '
# Row 37
$ws.Range("E37").Value = '  -2.57%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.61'
$ws.Range("E38").Value = '  +0.47%  '
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.835'
$ws.Range("E39").Value = '  -6.49%  '
# Row 40
$ws.Range("E40").Value = '  -2.74%  '
# Row 41
$ws.Range("E41").Value = '  -2.46%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '279.37'
$ws.Range("E42").Value = '  -5.06%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  +0.53%  '
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0983'
# Row 45
$ws.Range("E45").Value = '  -4.88%  '
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.48'
$ws.Range("E46").Value = '  -2.38%  '
# Row 47
$ws.Range("E47").Value = '  -4.98%  '
# Row 49
$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '1.984.16'
$ws.Range("E49").Value = '  +0.10%  '
# Row 50
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0227'
$ws.Range("E50").Value = '  -2.56%  '
# Row 51
$ws.Range("E51").Value = '  -4.10%  '
